$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shifted data rows (93-209)
$ws.Cells.Item(93, 4).Value = 44915
$ws.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(93, 9).Value = "Extra"
$ws.Cells.Item(93, 10).Value = 3000
$ws.Cells.Item(93, 11).Value = 2800
$ws.Cells.Item(93, 12).Value = 2800
$ws.Cells.Item(93, 13).Value = 2800
$ws.Cells.Item(93, 14).Value = "`$/unidad"
$ws.Cells.Item(93, 16).Value = 2800
$ws.Cells.Item(94, 4).Value = 44915
$ws.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 3000
$ws.Cells.Item(94, 11).Value = 2300
$ws.Cells.Item(94, 12).Value = 2300
$ws.Cells.Item(94, 13).Value = 2300
$ws.Cells.Item(94, 15).Value = "Región del Maule"
$ws.Cells.Item(94, 16).Value = 2300
$ws.Cells.Item(95, 4).Value = 44915
$ws.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(95, 9).Value = "Segunda"
$ws.Cells.Item(95, 10).Value = 2000
$ws.Cells.Item(95, 11).Value = 1800
$ws.Cells.Item(95, 12).Value = 1800
$ws.Cells.Item(95, 13).Value = 1800
$ws.Cells.Item(95, 15).Value = "Región del Maule"
$ws.Cells.Item(95, 16).Value = 1800
$ws.Cells.Item(96, 4).Value = 44169
$ws.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 2000
$ws.Cells.Item(96, 11).Value = 400
$ws.Cells.Item(96, 12).Value = 400
$ws.Cells.Item(96, 13).Value = 400
$ws.Cells.Item(96, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(96, 16).Value = 400
$ws.Cells.Item(97, 4).Value = 44901
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(97, 9).Value = "Especial"
$ws.Cells.Item(97, 10).Value = 2500
$ws.Cells.Item(97, 11).Value = 3200
$ws.Cells.Item(97, 12).Value = 3200
$ws.Cells.Item(97, 13).Value = 3200
$ws.Cells.Item(97, 15).Value = "Paine"
$ws.Cells.Item(97, 16).Value = 3200
$ws.Cells.Item(98, 4).Value = 44901
$ws.Cells.Item(98, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 2500
$ws.Cells.Item(98, 11).Value = 2800
$ws.Cells.Item(98, 12).Value = 2800
$ws.Cells.Item(98, 13).Value = 2800
$ws.Cells.Item(98, 15).Value = "Paine"
$ws.Cells.Item(98, 16).Value = 2800
$ws.Cells.Item(99, 4).Value = 44203
$ws.Cells.Item(99, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(99, 11).Value = 2200
$ws.Cells.Item(99, 12).Value = 2200
$ws.Cells.Item(99, 13).Value = 2200
$ws.Cells.Item(99, 16).Value = 2200
$ws.Cells.Item(100, 4).Value = 44203
$ws.Cells.Item(100, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(100, 10).Value = 5000
$ws.Cells.Item(100, 11).Value = 1800
$ws.Cells.Item(100, 12).Value = 1800
$ws.Cells.Item(100, 13).Value = 1800
$ws.Cells.Item(100, 16).Value = 1800
$ws.Cells.Item(101, 4).Value = 44203
$ws.Cells.Item(101, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(101, 11).Value = 1300
$ws.Cells.Item(101, 12).Value = 1300
$ws.Cells.Item(101, 13).Value = 1300
$ws.Cells.Item(101, 16).Value = 1300
$ws.Cells.Item(102, 4).Value = 44558
$ws.Cells.Item(102, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(102, 10).Value = 3000
$ws.Cells.Item(102, 11).Value = 2500
$ws.Cells.Item(102, 12).Value = 2500
$ws.Cells.Item(102, 13).Value = 2500
$ws.Cells.Item(102, 16).Value = 2500
$ws.Cells.Item(103, 4).Value = 44558
$ws.Cells.Item(103, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(103, 10).Value = 4000
$ws.Cells.Item(103, 11).Value = 2200
$ws.Cells.Item(103, 12).Value = 2200
$ws.Cells.Item(103, 13).Value = 2200
$ws.Cells.Item(103, 16).Value = 2200
$ws.Cells.Item(104, 4).Value = 44558
$ws.Cells.Item(104, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(104, 10).Value = 4000
$ws.Cells.Item(104, 11).Value = 1700
$ws.Cells.Item(104, 12).Value = 1700
$ws.Cells.Item(104, 13).Value = 1700
$ws.Cells.Item(104, 16).Value = 1700
$ws.Cells.Item(105, 4).Value = 44235
$ws.Cells.Item(105, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(105, 10).Value = 2000
$ws.Cells.Item(105, 11).Value = 2800
$ws.Cells.Item(105, 12).Value = 2800
$ws.Cells.Item(105, 13).Value = 2800
$ws.Cells.Item(105, 16).Value = 2800
$ws.Cells.Item(106, 4).Value = 44235
$ws.Cells.Item(106, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(106, 10).Value = 3000
$ws.Cells.Item(106, 11).Value = 2300
$ws.Cells.Item(106, 12).Value = 2300
$ws.Cells.Item(106, 13).Value = 2300
$ws.Cells.Item(106, 16).Value = 2300
$ws.Cells.Item(107, 4).Value = 44235
$ws.Cells.Item(107, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(107, 10).Value = 2000
$ws.Cells.Item(107, 11).Value = 1800
$ws.Cells.Item(107, 12).Value = 1800
$ws.Cells.Item(107, 13).Value = 1800
$ws.Cells.Item(107, 16).Value = 1800
$ws.Cells.Item(108, 4).Value = 44210
$ws.Cells.Item(108, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(108, 10).Value = 3000
$ws.Cells.Item(108, 11).Value = 2000
$ws.Cells.Item(108, 12).Value = 2000
$ws.Cells.Item(108, 13).Value = 2000
$ws.Cells.Item(108, 15).Value = "Región del Maule"
$ws.Cells.Item(108, 16).Value = 2000
$ws.Cells.Item(109, 4).Value = 44210
$ws.Cells.Item(109, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 6000
$ws.Cells.Item(109, 11).Value = 1600
$ws.Cells.Item(109, 12).Value = 1600
$ws.Cells.Item(109, 13).Value = 1600
$ws.Cells.Item(109, 16).Value = 1600
$ws.Cells.Item(110, 4).Value = 44210
$ws.Cells.Item(110, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(110, 9).Value = "Segunda"
$ws.Cells.Item(110, 10).Value = 4000
$ws.Cells.Item(110, 11).Value = 1200
$ws.Cells.Item(110, 12).Value = 1200
$ws.Cells.Item(110, 13).Value = 1200
$ws.Cells.Item(110, 15).Value = "Región del Maule"
$ws.Cells.Item(110, 16).Value = 1200
$ws.Cells.Item(111, 9).Value = "Extra"
$ws.Cells.Item(111, 10).Value = 2000
$ws.Cells.Item(111, 11).Value = 3000
$ws.Cells.Item(111, 12).Value = 3000
$ws.Cells.Item(111, 13).Value = 3000
$ws.Cells.Item(111, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(111, 16).Value = 3000
$ws.Cells.Item(112, 4).Value = 44546
$ws.Cells.Item(112, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(112, 10).Value = 2500
$ws.Cells.Item(112, 11).Value = 3000
$ws.Cells.Item(112, 12).Value = 3000
$ws.Cells.Item(112, 13).Value = 3000
$ws.Cells.Item(112, 16).Value = 3000
$ws.Cells.Item(113, 4).Value = 44546
$ws.Cells.Item(113, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(113, 11).Value = 2500
$ws.Cells.Item(113, 12).Value = 2500
$ws.Cells.Item(113, 13).Value = 2500
$ws.Cells.Item(113, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(113, 16).Value = 2500
$ws.Cells.Item(114, 4).Value = 44546
$ws.Cells.Item(114, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 3500
$ws.Cells.Item(114, 11).Value = 2500
$ws.Cells.Item(114, 12).Value = 2500
$ws.Cells.Item(114, 13).Value = 2500
$ws.Cells.Item(114, 16).Value = 2500
$ws.Cells.Item(115, 4).Value = 44568
$ws.Cells.Item(115, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(115, 10).Value = 2000
$ws.Cells.Item(115, 11).Value = 2300
$ws.Cells.Item(115, 12).Value = 2300
$ws.Cells.Item(115, 13).Value = 2300
$ws.Cells.Item(115, 16).Value = 2300
$ws.Cells.Item(116, 4).Value = 44568
$ws.Cells.Item(116, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(116, 10).Value = 3000
$ws.Cells.Item(116, 11).Value = 1800
$ws.Cells.Item(116, 12).Value = 1800
$ws.Cells.Item(116, 13).Value = 1800
$ws.Cells.Item(116, 16).Value = 1800
$ws.Cells.Item(117, 4).Value = 44568
$ws.Cells.Item(117, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(117, 9).Value = "Segunda"
$ws.Cells.Item(117, 11).Value = 1300
$ws.Cells.Item(117, 12).Value = 1300
$ws.Cells.Item(117, 13).Value = 1300
$ws.Cells.Item(117, 16).Value = 1300
$ws.Cells.Item(118, 4).Value = 44186
$ws.Cells.Item(118, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(118, 9).Value = "Extra"
$ws.Cells.Item(118, 10).Value = 2500
$ws.Cells.Item(118, 11).Value = 3200
$ws.Cells.Item(118, 12).Value = 3200
$ws.Cells.Item(118, 13).Value = 3200
$ws.Cells.Item(118, 16).Value = 3200
$ws.Cells.Item(119, 4).Value = 44186
$ws.Cells.Item(119, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 3500
$ws.Cells.Item(119, 11).Value = 2800
$ws.Cells.Item(119, 12).Value = 2800
$ws.Cells.Item(119, 13).Value = 2800
$ws.Cells.Item(119, 16).Value = 2800
$ws.Cells.Item(120, 4).Value = 44195
$ws.Cells.Item(120, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(120, 11).Value = 3000
$ws.Cells.Item(120, 12).Value = 3000
$ws.Cells.Item(120, 13).Value = 3000
$ws.Cells.Item(120, 16).Value = 3000
$ws.Cells.Item(121, 4).Value = 44195
$ws.Cells.Item(121, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(121, 11).Value = 2500
$ws.Cells.Item(121, 12).Value = 2500
$ws.Cells.Item(121, 13).Value = 2500
$ws.Cells.Item(121, 16).Value = 2500
$ws.Cells.Item(122, 4).Value = 44195
$ws.Cells.Item(122, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(122, 11).Value = 2000
$ws.Cells.Item(122, 12).Value = 2000
$ws.Cells.Item(122, 13).Value = 2000
$ws.Cells.Item(122, 16).Value = 2000
$ws.Cells.Item(123, 4).Value = 44567
$ws.Cells.Item(123, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(123, 10).Value = 2000
$ws.Cells.Item(123, 11).Value = 2300
$ws.Cells.Item(123, 12).Value = 2300
$ws.Cells.Item(123, 13).Value = 2300
$ws.Cells.Item(123, 16).Value = 2300
$ws.Cells.Item(124, 4).Value = 44567
$ws.Cells.Item(124, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 10).Value = 3000
$ws.Cells.Item(124, 11).Value = 1800
$ws.Cells.Item(124, 12).Value = 1800
$ws.Cells.Item(124, 13).Value = 1800
$ws.Cells.Item(124, 16).Value = 1800
$ws.Cells.Item(125, 4).Value = 44567
$ws.Cells.Item(125, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 10).Value = 2000
$ws.Cells.Item(125, 11).Value = 1300
$ws.Cells.Item(125, 12).Value = 1300
$ws.Cells.Item(125, 13).Value = 1300
$ws.Cells.Item(125, 16).Value = 1300
$ws.Cells.Item(126, 4).Value = 44214
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 2000
$ws.Cells.Item(126, 12).Value = 2000
$ws.Cells.Item(126, 13).Value = 2000
$ws.Cells.Item(126, 16).Value = 2000
$ws.Cells.Item(127, 4).Value = 44214
$ws.Cells.Item(127, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(127, 10).Value = 8000
$ws.Cells.Item(127, 11).Value = 1600
$ws.Cells.Item(127, 12).Value = 1600
$ws.Cells.Item(127, 13).Value = 1600
$ws.Cells.Item(127, 16).Value = 1600
$ws.Cells.Item(128, 4).Value = 44214
$ws.Cells.Item(128, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(128, 10).Value = 4000
$ws.Cells.Item(128, 11).Value = 1200
$ws.Cells.Item(128, 12).Value = 1200
$ws.Cells.Item(128, 13).Value = 1200
$ws.Cells.Item(128, 16).Value = 1200
$ws.Cells.Item(129, 4).Value = 44202
$ws.Cells.Item(129, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(129, 10).Value = 5000
$ws.Cells.Item(130, 4).Value = 44202
$ws.Cells.Item(130, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(130, 10).Value = 10000
$ws.Cells.Item(130, 11).Value = 2000
$ws.Cells.Item(130, 12).Value = 2000
$ws.Cells.Item(130, 13).Value = 2000
$ws.Cells.Item(130, 16).Value = 2000
$ws.Cells.Item(131, 4).Value = 44202
$ws.Cells.Item(131, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(131, 10).Value = 8000
$ws.Cells.Item(131, 11).Value = 1800
$ws.Cells.Item(131, 12).Value = 1800
$ws.Cells.Item(131, 13).Value = 1800
$ws.Cells.Item(131, 16).Value = 1800
$ws.Cells.Item(132, 4).Value = 44572
$ws.Cells.Item(132, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 2300
$ws.Cells.Item(132, 12).Value = 2300
$ws.Cells.Item(132, 13).Value = 2300
$ws.Cells.Item(132, 16).Value = 2300
$ws.Cells.Item(133, 4).Value = 44572
$ws.Cells.Item(133, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(133, 10).Value = 4000
$ws.Cells.Item(133, 11).Value = 1800
$ws.Cells.Item(133, 12).Value = 1800
$ws.Cells.Item(133, 13).Value = 1800
$ws.Cells.Item(133, 16).Value = 1800
$ws.Cells.Item(134, 4).Value = 44572
$ws.Cells.Item(134, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 1500
$ws.Cells.Item(134, 12).Value = 1500
$ws.Cells.Item(134, 13).Value = 1500
$ws.Cells.Item(134, 16).Value = 1500
$ws.Cells.Item(135, 4).Value = 44209
$ws.Cells.Item(135, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(135, 11).Value = 2000
$ws.Cells.Item(135, 12).Value = 2000
$ws.Cells.Item(135, 13).Value = 2000
$ws.Cells.Item(135, 16).Value = 2000
$ws.Cells.Item(136, 4).Value = 44209
$ws.Cells.Item(136, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 1600
$ws.Cells.Item(136, 12).Value = 1600
$ws.Cells.Item(136, 13).Value = 1600
$ws.Cells.Item(136, 16).Value = 1600
$ws.Cells.Item(137, 4).Value = 44209
$ws.Cells.Item(137, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(137, 9).Value = "Segunda"
$ws.Cells.Item(137, 10).Value = 5000
$ws.Cells.Item(137, 11).Value = 1200
$ws.Cells.Item(137, 12).Value = 1200
$ws.Cells.Item(137, 13).Value = 1200
$ws.Cells.Item(137, 16).Value = 1200
$ws.Cells.Item(138, 4).Value = 44554
$ws.Cells.Item(138, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(138, 9).Value = "Extra"
$ws.Cells.Item(138, 10).Value = 3000
$ws.Cells.Item(138, 11).Value = 2500
$ws.Cells.Item(138, 12).Value = 2500
$ws.Cells.Item(138, 13).Value = 2500
$ws.Cells.Item(138, 16).Value = 2500
$ws.Cells.Item(139, 4).Value = 44554
$ws.Cells.Item(139, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 3000
$ws.Cells.Item(139, 11).Value = 1900
$ws.Cells.Item(139, 12).Value = 1900
$ws.Cells.Item(139, 13).Value = 1900
$ws.Cells.Item(139, 16).Value = 1900
$ws.Cells.Item(140, 4).Value = 44200
$ws.Cells.Item(140, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(140, 10).Value = 4000
$ws.Cells.Item(140, 11).Value = 2200
$ws.Cells.Item(140, 12).Value = 2200
$ws.Cells.Item(140, 13).Value = 2200
$ws.Cells.Item(140, 16).Value = 2200
$ws.Cells.Item(141, 4).Value = 44200
$ws.Cells.Item(141, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(141, 10).Value = 8000
$ws.Cells.Item(141, 11).Value = 1700
$ws.Cells.Item(141, 12).Value = 1700
$ws.Cells.Item(141, 13).Value = 1700
$ws.Cells.Item(141, 16).Value = 1700
$ws.Cells.Item(142, 4).Value = 44200
$ws.Cells.Item(142, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(142, 10).Value = 4000
$ws.Cells.Item(142, 11).Value = 1400
$ws.Cells.Item(142, 12).Value = 1400
$ws.Cells.Item(142, 13).Value = 1400
$ws.Cells.Item(142, 16).Value = 1400
$ws.Cells.Item(143, 4).Value = 44573
$ws.Cells.Item(143, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(143, 11).Value = 2300
$ws.Cells.Item(143, 12).Value = 2300
$ws.Cells.Item(143, 13).Value = 2300
$ws.Cells.Item(143, 16).Value = 2300
$ws.Cells.Item(144, 4).Value = 44573
$ws.Cells.Item(144, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(144, 11).Value = 1800
$ws.Cells.Item(144, 12).Value = 1800
$ws.Cells.Item(144, 13).Value = 1800
$ws.Cells.Item(144, 16).Value = 1800
$ws.Cells.Item(145, 4).Value = 44573
$ws.Cells.Item(145, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(146, 4).Value = 44560
$ws.Cells.Item(146, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(146, 10).Value = 2000
$ws.Cells.Item(147, 4).Value = 44560
$ws.Cells.Item(147, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(147, 10).Value = 3000
$ws.Cells.Item(147, 11).Value = 2000
$ws.Cells.Item(147, 12).Value = 2000
$ws.Cells.Item(147, 13).Value = 2000
$ws.Cells.Item(147, 16).Value = 2000
$ws.Cells.Item(148, 4).Value = 44560
$ws.Cells.Item(148, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(148, 10).Value = 2000
$ws.Cells.Item(148, 11).Value = 1500
$ws.Cells.Item(148, 12).Value = 1500
$ws.Cells.Item(148, 13).Value = 1500
$ws.Cells.Item(148, 16).Value = 1500
$ws.Cells.Item(149, 4).Value = 44579
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 10).Value = 3000
$ws.Cells.Item(149, 11).Value = 2500
$ws.Cells.Item(149, 12).Value = 2500
$ws.Cells.Item(149, 13).Value = 2500
$ws.Cells.Item(149, 16).Value = 2500
$ws.Cells.Item(150, 4).Value = 44579
$ws.Cells.Item(150, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(150, 10).Value = 2100
$ws.Cells.Item(151, 4).Value = 44579
$ws.Cells.Item(151, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(151, 10).Value = 1300
$ws.Cells.Item(152, 4).Value = 44585
$ws.Cells.Item(152, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(152, 11).Value = 2000
$ws.Cells.Item(152, 12).Value = 2000
$ws.Cells.Item(152, 13).Value = 2000
$ws.Cells.Item(152, 16).Value = 2000
$ws.Cells.Item(153, 4).Value = 44585
$ws.Cells.Item(153, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(153, 10).Value = 3000
$ws.Cells.Item(153, 11).Value = 1500
$ws.Cells.Item(153, 12).Value = 1500
$ws.Cells.Item(153, 13).Value = 1500
$ws.Cells.Item(153, 16).Value = 1500
$ws.Cells.Item(154, 4).Value = 44585
$ws.Cells.Item(154, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(154, 10).Value = 2000
$ws.Cells.Item(154, 11).Value = 1000
$ws.Cells.Item(154, 12).Value = 1000
$ws.Cells.Item(154, 13).Value = 1000
$ws.Cells.Item(154, 16).Value = 1000
$ws.Cells.Item(155, 4).Value = 44232
$ws.Cells.Item(155, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(155, 9).Value = "Extra"
$ws.Cells.Item(155, 11).Value = 2500
$ws.Cells.Item(155, 12).Value = 2500
$ws.Cells.Item(155, 13).Value = 2500
$ws.Cells.Item(155, 16).Value = 2500
$ws.Cells.Item(156, 4).Value = 44232
$ws.Cells.Item(156, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(156, 10).Value = 5000
$ws.Cells.Item(156, 11).Value = 2000
$ws.Cells.Item(156, 12).Value = 2000
$ws.Cells.Item(156, 13).Value = 2000
$ws.Cells.Item(156, 16).Value = 2000
$ws.Cells.Item(157, 4).Value = 44232
$ws.Cells.Item(157, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(157, 9).Value = "Segunda"
$ws.Cells.Item(157, 11).Value = 1600
$ws.Cells.Item(157, 12).Value = 1600
$ws.Cells.Item(157, 13).Value = 1600
$ws.Cells.Item(157, 16).Value = 1600
$ws.Cells.Item(158, 4).Value = 44179
$ws.Cells.Item(158, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(158, 9).Value = "Especial"
$ws.Cells.Item(158, 10).Value = 2000
$ws.Cells.Item(158, 11).Value = 3200
$ws.Cells.Item(158, 12).Value = 3200
$ws.Cells.Item(158, 13).Value = 3200
$ws.Cells.Item(158, 16).Value = 3200
$ws.Cells.Item(159, 4).Value = 44179
$ws.Cells.Item(159, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 11).Value = 2500
$ws.Cells.Item(159, 12).Value = 2500
$ws.Cells.Item(159, 13).Value = 2500
$ws.Cells.Item(159, 16).Value = 2500
$ws.Cells.Item(160, 4).Value = 44221
$ws.Cells.Item(160, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(160, 10).Value = 3000
$ws.Cells.Item(160, 11).Value = 2300
$ws.Cells.Item(160, 12).Value = 2300
$ws.Cells.Item(160, 13).Value = 2300
$ws.Cells.Item(160, 16).Value = 2300
$ws.Cells.Item(161, 4).Value = 44221
$ws.Cells.Item(161, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(161, 10).Value = 5000
$ws.Cells.Item(161, 11).Value = 1800
$ws.Cells.Item(161, 12).Value = 1800
$ws.Cells.Item(161, 13).Value = 1800
$ws.Cells.Item(161, 16).Value = 1800
$ws.Cells.Item(162, 4).Value = 44221
$ws.Cells.Item(162, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(162, 11).Value = 1300
$ws.Cells.Item(162, 12).Value = 1300
$ws.Cells.Item(162, 13).Value = 1300
$ws.Cells.Item(162, 16).Value = 1300
$ws.Cells.Item(163, 4).Value = 44599
$ws.Cells.Item(163, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(163, 11).Value = 2000
$ws.Cells.Item(163, 12).Value = 2000
$ws.Cells.Item(163, 13).Value = 2000
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 2000
$ws.Cells.Item(164, 4).Value = 44599
$ws.Cells.Item(164, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(164, 10).Value = 3000
$ws.Cells.Item(164, 11).Value = 1500
$ws.Cells.Item(164, 12).Value = 1500
$ws.Cells.Item(164, 13).Value = 1500
$ws.Cells.Item(164, 15).Value = "Región del Maule"
$ws.Cells.Item(164, 16).Value = 1500
$ws.Cells.Item(165, 4).Value = 44599
$ws.Cells.Item(165, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(165, 10).Value = 3000
$ws.Cells.Item(165, 11).Value = 1000
$ws.Cells.Item(165, 12).Value = 1000
$ws.Cells.Item(165, 13).Value = 1000
$ws.Cells.Item(165, 15).Value = "Región del Maule"
$ws.Cells.Item(165, 16).Value = 1000
$ws.Cells.Item(166, 4).Value = 44908
$ws.Cells.Item(166, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(166, 11).Value = 3500
$ws.Cells.Item(166, 12).Value = 3500
$ws.Cells.Item(166, 13).Value = 3500
$ws.Cells.Item(166, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(166, 16).Value = 3500
$ws.Cells.Item(167, 4).Value = 44908
$ws.Cells.Item(167, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(167, 10).Value = 2000
$ws.Cells.Item(167, 11).Value = 3000
$ws.Cells.Item(167, 12).Value = 3000
$ws.Cells.Item(167, 13).Value = 3000
$ws.Cells.Item(167, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(167, 16).Value = 3000
$ws.Cells.Item(168, 4).Value = 44908
$ws.Cells.Item(168, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(168, 11).Value = 2500
$ws.Cells.Item(168, 12).Value = 2500
$ws.Cells.Item(168, 13).Value = 2500
$ws.Cells.Item(168, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(168, 16).Value = 2500
$ws.Cells.Item(169, 4).Value = 44194
$ws.Cells.Item(169, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(169, 9).Value = "Extra"
$ws.Cells.Item(169, 10).Value = 2000
$ws.Cells.Item(170, 4).Value = 44194
$ws.Cells.Item(170, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 4000
$ws.Cells.Item(171, 4).Value = 44194
$ws.Cells.Item(171, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(171, 9).Value = "Segunda"
$ws.Cells.Item(171, 10).Value = 2000
$ws.Cells.Item(171, 11).Value = 2000
$ws.Cells.Item(171, 12).Value = 2000
$ws.Cells.Item(171, 13).Value = 2000
$ws.Cells.Item(171, 15).Value = "Región del Maule"
$ws.Cells.Item(171, 16).Value = 2000
$ws.Cells.Item(172, 4).Value = 44181
$ws.Cells.Item(172, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 5000
$ws.Cells.Item(172, 11).Value = 3000
$ws.Cells.Item(172, 12).Value = 3000
$ws.Cells.Item(172, 13).Value = 3000
$ws.Cells.Item(172, 15).Value = "Región del Maule"
$ws.Cells.Item(172, 16).Value = 3000
$ws.Cells.Item(173, 4).Value = 44181
$ws.Cells.Item(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(173, 9).Value = "Segunda"
$ws.Cells.Item(173, 10).Value = 3000
$ws.Cells.Item(173, 11).Value = 2500
$ws.Cells.Item(173, 12).Value = 2500
$ws.Cells.Item(173, 13).Value = 2500
$ws.Cells.Item(173, 16).Value = 2500
$ws.Cells.Item(174, 4).Value = 44907
$ws.Cells.Item(174, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(174, 10).Value = 2500
$ws.Cells.Item(174, 11).Value = 3000
$ws.Cells.Item(174, 12).Value = 3000
$ws.Cells.Item(174, 13).Value = 3000
$ws.Cells.Item(174, 15).Value = "Paine"
$ws.Cells.Item(174, 16).Value = 3000
$ws.Cells.Item(175, 4).Value = 44907
$ws.Cells.Item(175, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(175, 9).Value = "Segunda"
$ws.Cells.Item(175, 10).Value = 2500
$ws.Cells.Item(175, 11).Value = 2500
$ws.Cells.Item(175, 12).Value = 2500
$ws.Cells.Item(175, 13).Value = 2500
$ws.Cells.Item(175, 15).Value = "Paine"
$ws.Cells.Item(175, 16).Value = 2500
$ws.Cells.Item(176, 4).Value = 44578
$ws.Cells.Item(176, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(176, 9).Value = "Extra"
$ws.Cells.Item(176, 10).Value = 2500
$ws.Cells.Item(176, 11).Value = 2000
$ws.Cells.Item(176, 12).Value = 2000
$ws.Cells.Item(176, 13).Value = 2000
$ws.Cells.Item(176, 16).Value = 2000
$ws.Cells.Item(177, 4).Value = 44578
$ws.Cells.Item(177, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 2000
$ws.Cells.Item(177, 11).Value = 1500
$ws.Cells.Item(177, 12).Value = 1500
$ws.Cells.Item(177, 13).Value = 1500
$ws.Cells.Item(177, 16).Value = 1500
$ws.Cells.Item(178, 4).Value = 44566
$ws.Cells.Item(178, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(178, 9).Value = "Extra"
$ws.Cells.Item(178, 11).Value = 2000
$ws.Cells.Item(178, 12).Value = 2000
$ws.Cells.Item(178, 13).Value = 2000
$ws.Cells.Item(178, 14).Value = "`$/unidad"
$ws.Cells.Item(178, 15).Value = "Región del Maule"
$ws.Cells.Item(178, 16).Value = 2000
$ws.Cells.Item(179, 4).Value = 44566
$ws.Cells.Item(179, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(179, 10).Value = 3000
$ws.Cells.Item(179, 11).Value = 1500
$ws.Cells.Item(179, 12).Value = 1500
$ws.Cells.Item(179, 13).Value = 1500
$ws.Cells.Item(179, 14).Value = "`$/unidad"
$ws.Cells.Item(179, 15).Value = "Región del Maule"
$ws.Cells.Item(179, 16).Value = 1500
$ws.Cells.Item(180, 4).Value = 44566
$ws.Cells.Item(180, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(180, 9).Value = "Segunda"
$ws.Cells.Item(180, 10).Value = 3000
$ws.Cells.Item(180, 11).Value = 1000
$ws.Cells.Item(180, 12).Value = 1000
$ws.Cells.Item(180, 13).Value = 1000
$ws.Cells.Item(180, 16).Value = 1000
$ws.Cells.Item(181, 4).Value = 44895
$ws.Cells.Item(181, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(181, 10).Value = 2000
$ws.Cells.Item(181, 11).Value = 400
$ws.Cells.Item(181, 12).Value = 400
$ws.Cells.Item(181, 13).Value = 400
$ws.Cells.Item(181, 14).Value = "`$/kilo"
$ws.Cells.Item(181, 15).Value = "Paine"
$ws.Cells.Item(181, 16).Value = 400
$ws.Cells.Item(182, 4).Value = 44895
$ws.Cells.Item(182, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 2000
$ws.Cells.Item(182, 11).Value = 500
$ws.Cells.Item(182, 12).Value = 500
$ws.Cells.Item(182, 13).Value = 500
$ws.Cells.Item(182, 14).Value = "`$/kilo"
$ws.Cells.Item(182, 15).Value = "Perú"
$ws.Cells.Item(182, 16).Value = 500
$ws.Cells.Item(183, 4).Value = 44211
$ws.Cells.Item(183, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(183, 10).Value = 5000
$ws.Cells.Item(183, 11).Value = 2000
$ws.Cells.Item(183, 12).Value = 2000
$ws.Cells.Item(183, 13).Value = 2000
$ws.Cells.Item(183, 16).Value = 2000
$ws.Cells.Item(184, 4).Value = 44211
$ws.Cells.Item(184, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(184, 10).Value = 8000
$ws.Cells.Item(184, 11).Value = 1600
$ws.Cells.Item(184, 12).Value = 1600
$ws.Cells.Item(184, 13).Value = 1600
$ws.Cells.Item(184, 16).Value = 1600
$ws.Cells.Item(185, 4).Value = 44211
$ws.Cells.Item(185, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(185, 10).Value = 5000
$ws.Cells.Item(185, 11).Value = 1200
$ws.Cells.Item(185, 12).Value = 1200
$ws.Cells.Item(185, 13).Value = 1200
$ws.Cells.Item(185, 16).Value = 1200
$ws.Cells.Item(186, 4).Value = 44559
$ws.Cells.Item(186, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(186, 11).Value = 2300
$ws.Cells.Item(186, 12).Value = 2300
$ws.Cells.Item(186, 13).Value = 2300
$ws.Cells.Item(186, 16).Value = 2300
$ws.Cells.Item(187, 4).Value = 44559
$ws.Cells.Item(187, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(187, 10).Value = 3000
$ws.Cells.Item(187, 11).Value = 2000
$ws.Cells.Item(187, 12).Value = 2000
$ws.Cells.Item(187, 13).Value = 2000
$ws.Cells.Item(187, 16).Value = 2000
$ws.Cells.Item(188, 4).Value = 44559
$ws.Cells.Item(188, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(188, 11).Value = 1500
$ws.Cells.Item(188, 12).Value = 1500
$ws.Cells.Item(188, 13).Value = 1500
$ws.Cells.Item(188, 16).Value = 1500
$ws.Cells.Item(189, 4).Value = 44216
$ws.Cells.Item(189, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(189, 10).Value = 3000
$ws.Cells.Item(189, 11).Value = 1800
$ws.Cells.Item(189, 12).Value = 1800
$ws.Cells.Item(189, 13).Value = 1800
$ws.Cells.Item(189, 16).Value = 1800
$ws.Cells.Item(190, 4).Value = 44216
$ws.Cells.Item(190, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(190, 10).Value = 5000
$ws.Cells.Item(190, 11).Value = 1400
$ws.Cells.Item(190, 12).Value = 1400
$ws.Cells.Item(190, 13).Value = 1400
$ws.Cells.Item(190, 16).Value = 1400
$ws.Cells.Item(191, 4).Value = 44216
$ws.Cells.Item(191, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(191, 10).Value = 3000
$ws.Cells.Item(191, 11).Value = 1000
$ws.Cells.Item(191, 12).Value = 1000
$ws.Cells.Item(191, 13).Value = 1000
$ws.Cells.Item(191, 16).Value = 1000
$ws.Cells.Item(192, 4).Value = 44264
$ws.Cells.Item(192, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(192, 11).Value = 2300
$ws.Cells.Item(192, 12).Value = 2300
$ws.Cells.Item(192, 13).Value = 2300
$ws.Cells.Item(192, 16).Value = 2300
$ws.Cells.Item(193, 4).Value = 44264
$ws.Cells.Item(193, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(193, 10).Value = 2000
$ws.Cells.Item(193, 11).Value = 1800
$ws.Cells.Item(193, 12).Value = 1800
$ws.Cells.Item(193, 13).Value = 1800
$ws.Cells.Item(193, 16).Value = 1800
$ws.Cells.Item(194, 4).Value = 44264
$ws.Cells.Item(194, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(194, 10).Value = 2000
$ws.Cells.Item(194, 11).Value = 1300
$ws.Cells.Item(194, 12).Value = 1300
$ws.Cells.Item(194, 13).Value = 1300
$ws.Cells.Item(194, 16).Value = 1300
$ws.Cells.Item(195, 4).Value = 44592
$ws.Cells.Item(195, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(195, 9).Value = "Extra"
$ws.Cells.Item(195, 11).Value = 2000
$ws.Cells.Item(195, 12).Value = 2000
$ws.Cells.Item(195, 13).Value = 2000
$ws.Cells.Item(195, 14).Value = "`$/unidad"
$ws.Cells.Item(195, 15).Value = "Región del Maule"
$ws.Cells.Item(195, 16).Value = 2000
$ws.Cells.Item(196, 4).Value = 44592
$ws.Cells.Item(196, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 11).Value = 1500
$ws.Cells.Item(196, 12).Value = 1500
$ws.Cells.Item(196, 13).Value = 1500
$ws.Cells.Item(196, 16).Value = 1500
$ws.Cells.Item(197, 4).Value = 44592
$ws.Cells.Item(197, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(197, 9).Value = "Segunda"
$ws.Cells.Item(197, 10).Value = 3000
$ws.Cells.Item(197, 11).Value = 1000
$ws.Cells.Item(197, 12).Value = 1000
$ws.Cells.Item(197, 13).Value = 1000
$ws.Cells.Item(197, 16).Value = 1000
$ws.Cells.Item(198, 4).Value = 44533
$ws.Cells.Item(198, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 2000
$ws.Cells.Item(198, 11).Value = 500
$ws.Cells.Item(198, 12).Value = 500
$ws.Cells.Item(198, 13).Value = 500
$ws.Cells.Item(198, 14).Value = "`$/kilo"
$ws.Cells.Item(198, 15).Value = "Perú"
$ws.Cells.Item(198, 16).Value = 500
$ws.Cells.Item(199, 4).Value = 44217
$ws.Cells.Item(199, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(199, 10).Value = 3000
$ws.Cells.Item(200, 4).Value = 44217
$ws.Cells.Item(200, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(200, 10).Value = 8000
$ws.Cells.Item(200, 11).Value = 1600
$ws.Cells.Item(200, 12).Value = 1600
$ws.Cells.Item(200, 13).Value = 1600
$ws.Cells.Item(200, 16).Value = 1600
$ws.Cells.Item(201, 4).Value = 44217
$ws.Cells.Item(201, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(201, 9).Value = "Segunda"
$ws.Cells.Item(201, 10).Value = 4000
$ws.Cells.Item(201, 11).Value = 1300
$ws.Cells.Item(201, 12).Value = 1300
$ws.Cells.Item(201, 13).Value = 1300
$ws.Cells.Item(201, 16).Value = 1300
$ws.Cells.Item(202, 4).Value = 44580
$ws.Cells.Item(202, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(202, 9).Value = "Extra"
$ws.Cells.Item(202, 10).Value = 2500
$ws.Cells.Item(203, 4).Value = 44580
$ws.Cells.Item(203, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 3500
$ws.Cells.Item(204, 4).Value = 44565
$ws.Cells.Item(204, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(204, 11).Value = 2500
$ws.Cells.Item(204, 12).Value = 2500
$ws.Cells.Item(204, 13).Value = 2500
$ws.Cells.Item(204, 16).Value = 2500
$ws.Cells.Item(205, 4).Value = 44565
$ws.Cells.Item(205, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(205, 11).Value = 2000
$ws.Cells.Item(205, 12).Value = 2000
$ws.Cells.Item(205, 13).Value = 2000
$ws.Cells.Item(205, 16).Value = 2000
$ws.Cells.Item(206, 4).Value = 44565
$ws.Cells.Item(206, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(207, 4).Value = 44571
$ws.Cells.Item(207, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(207, 10).Value = 2000
$ws.Cells.Item(207, 11).Value = 2300
$ws.Cells.Item(207, 12).Value = 2300
$ws.Cells.Item(207, 13).Value = 2300
$ws.Cells.Item(207, 15).Value = "Región del Maule"
$ws.Cells.Item(207, 16).Value = 2300
$ws.Cells.Item(208, 4).Value = 44571
$ws.Cells.Item(208, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(208, 10).Value = 3000
$ws.Cells.Item(208, 11).Value = 1800
$ws.Cells.Item(208, 12).Value = 1800
$ws.Cells.Item(208, 13).Value = 1800
$ws.Cells.Item(208, 15).Value = "Región del Maule"
$ws.Cells.Item(208, 16).Value = 1800
$ws.Cells.Item(209, 4).Value = 44571
$ws.Cells.Item(209, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(209, 10).Value = 2000
$ws.Cells.Item(209, 11).Value = 1500
$ws.Cells.Item(209, 12).Value = 1500
$ws.Cells.Item(209, 13).Value = 1500
$ws.Cells.Item(209, 15).Value = "Región del Maule"
$ws.Cells.Item(209, 16).Value = 1500

# Append new rows (210-212)
$ws.Cells.Item(210, 1).Value = 5
$ws.Cells.Item(210, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(210, 3).Value = "Maule"
$ws.Cells.Item(210, 4).Value = 44911
$ws.Cells.Item(210, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(210, 5).Value = 7
$ws.Cells.Item(210, 6).Value = 100112028
$ws.Cells.Item(210, 7).Value = "Sandia"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Extra"
$ws.Cells.Item(210, 10).Value = 5000
$ws.Cells.Item(210, 11).Value = 2800
$ws.Cells.Item(210, 12).Value = 2800
$ws.Cells.Item(210, 13).Value = 2800
$ws.Cells.Item(210, 14).Value = "`$/unidad"
$ws.Cells.Item(210, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(210, 16).Value = 2800
$ws.Cells.Item(210, 17).Value = 1
$ws.Cells.Item(210, 18).Value = "Hortaliza"
$ws.Cells.Item(211, 1).Value = 5
$ws.Cells.Item(211, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(211, 3).Value = "Maule"
$ws.Cells.Item(211, 4).Value = 44911
$ws.Cells.Item(211, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(211, 5).Value = 7
$ws.Cells.Item(211, 6).Value = 100112028
$ws.Cells.Item(211, 7).Value = "Sandia"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 5000
$ws.Cells.Item(211, 11).Value = 2300
$ws.Cells.Item(211, 12).Value = 2300
$ws.Cells.Item(211, 13).Value = 2300
$ws.Cells.Item(211, 14).Value = "`$/unidad"
$ws.Cells.Item(211, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(211, 16).Value = 2300
$ws.Cells.Item(211, 17).Value = 1
$ws.Cells.Item(211, 18).Value = "Hortaliza"
$ws.Cells.Item(212, 1).Value = 5
$ws.Cells.Item(212, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(212, 3).Value = "Maule"
$ws.Cells.Item(212, 4).Value = 44911
$ws.Cells.Item(212, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212, 5).Value = 7
$ws.Cells.Item(212, 6).Value = 100112028
$ws.Cells.Item(212, 7).Value = "Sandia"
$ws.Cells.Item(212, 8).Value = "Sin especificar"
$ws.Cells.Item(212, 9).Value = "Segunda"
$ws.Cells.Item(212, 10).Value = 5000
$ws.Cells.Item(212, 11).Value = 1800
$ws.Cells.Item(212, 12).Value = 1800
$ws.Cells.Item(212, 13).Value = 1800
$ws.Cells.Item(212, 14).Value = "`$/unidad"
$ws.Cells.Item(212, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(212, 16).Value = 1800
$ws.Cells.Item(212, 17).Value = 1
$ws.Cells.Item(212, 18).Value = "Hortaliza"
